$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.892.71"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "3.800.67"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'601.35"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'172.31"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("D7").Value = "3.798.20"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  -4.27%  "
$ws.Range("E11").Value = "  -5.92%  "
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("D13").Value = "'38.82"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "4.438.05"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "3.800.07"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "67.852.01"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").Value = "'7.27"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("E19").Value = "  -3.88%  "
$ws.Range("D20").Value = "'17.35"
$ws.Range("E20").Value = "  +6.09%  "
$ws.Range("D21").Value = "'494.53"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").Value = "'9.19"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "'0.742"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").Value = "'85.82"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  -4.90%  "
$ws.Range("E26").Value = "  +7.62%  "
$ws.Range("D27").Value = "'12.39"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D32").Value = "'32.95"
$ws.Range("E32").Value = "  +7.51%  "
$ws.Range("D33").Value = "'7.87"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "'0.110"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").Value = "'5.85"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'464.64"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.332"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("D42").Value = "'49.04"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").Value = "'8.44"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").Value = "'41.15"
$ws.Range("E45").Value = "  -8.04%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "2.846.87"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").Value = "'140.18"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "'0.0353"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").Value = "'25.91"
$ws.Range("E50").Value = "  -5.42%  "
$ws.Range("D51").Value = "'24.23"
$ws.Range("E51").Value = "  +11.52%  "
